$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header columns:
#   Before: A1=Email, B1=Games, C1=Name, D1=Age
#   After:  A1=Name,  B1=Age,   C1=Email, D1=Games
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Games"

# The "Games" list validation dropdown now lives in column D (it used to be
# in column B), so move the data validation range accordingly.
$ws.Range("B2:B100000").Validation.Delete()
$ws.Range("D2:D100000").Validation.Add(3, 1, 1, "Super Mario,SONIC,Zelda,GTA")
